$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E2").Value = "smallest (N=28)"
$ws.Range("E7").Value = "small (N=615)"

$ws.Range("D9").Select()
